$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

$rng = $ws1.Range("P60:Q69")
$rng.WrapText = $true
$rng.Borders.Color = 13421772
$rng.Borders.Weight = -4138
Write-Host "done"
